{"js": "// Remove the trailing \"Ver no Jupiter ...\" / \"\u00a9 2020 ...\" footer block\n// (and the blank paragraph right before it) that used to follow the\n// \"2014.\" paragraph at the end of the document body.\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"text\");\nawait context.sync();\n\nconst items = paragraphs.items;\n\n// Locate the paragraph whose text is exactly \"2014.\" \u2014 the anchor point\n// right before the block we need to drop.\nlet anchorIndex = -1;\nfor (let i = 0; i < items.length; i++) {\n  if (items[i].text === \"2014.\") {\n    anchorIndex = i;\n    break;\n  }\n}\n\nif (anchorIndex !== -1) {\n  // The three paragraphs right after \"2014.\" are: an empty paragraph,\n  // \"Ver no Jupiter Salvar em pdf Salvar em docx\", and the \"\u00a9 2020 ...\"\n  // copyright/footer line. Delete them (in reverse order so indices stay\n  // valid while deleting).\n  const toDelete = [];\n  for (let j = anchorIndex + 1; j < items.length; j++) {\n    const t = items[j].text;\n    if (\n      t === \"\" ||\n      t === \"Ver no Jupiter Salvar em pdf Salvar em docx\" ||\n      t ===\n        \"\u00a9 2020 . Contact: luizeleno@usp.br. Powered by Jekyll and Github pages. Original theme under Creative Commons Attribution\"\n    ) {\n      toDelete.push(items[j]);\n      // Stop once we've collected the blank line + the two text lines.\n      if (toDelete.length === 3) break;\n    } else {\n      break;\n    }\n  }\n\n  for (let k = toDelete.length - 1; k >= 0; k--) {\n    toDelete[k].delete();\n  }\n\n  await context.sync();\n}\n", "ps1": "# Remove the trailing \"Ver no Jupiter ...\" / \"C 2020 ...\" footer block\n# (and the blank paragraph right before it) that used to follow the\n# \"2014.\" paragraph at the end of the document body.\n$d = $word.ActiveDocument\n\n# Find the paragraph whose text is exactly \"2014.\" - the anchor point\n# right before the block we need to drop.\n$anchorIndex = -1\n$count = $d.Paragraphs.Count\nfor ($i = 1; $i -le $count; $i++) {\n    $p = $d.Paragraphs.Item($i)\n    if ($p.Range.Text.Trim() -eq \"2014.\") {\n        $anchorIndex = $i\n        break\n    }\n}\n\nif ($anchorIndex -ge 1) {\n    # The three paragraphs right after \"2014.\" are: an empty paragraph,\n    # \"Ver no Jupiter Salvar em pdf Salvar em docx\", and the \"(c) 2020 ...\"\n    # copyright/footer line. Collect them and delete as one range.\n    $ver = \"Ver no Jupiter Salvar em pdf Salvar em docx\"\n    $copyright = [char]0x00A9 + \" 2020 . Contact: luizeleno@usp.br. Powered by Jekyll and Github pages. Original theme under Creative Commons Attribution\"\n\n    $toDelete = @()\n    $j = $anchorIndex + 1\n    while ($j -le $d.Paragraphs.Count -and $toDelete.Count -lt 3) {\n        $t = $d.Paragraphs.Item($j).Range.Text.Trim()\n        if ($t -eq \"\" -or $t -eq $ver -or $t -eq $copyright) {\n            $toDelete += $j\n            $j = $j + 1\n        } else {\n            break\n        }\n    }\n\n    if ($toDelete.Count -gt 0) {\n        $firstPara = $d.Paragraphs.Item($toDelete[0])\n        $lastPara = $d.Paragraphs.Item($toDelete[$toDelete.Count - 1])\n        $delRange = $d.Range($firstPara.Range.Start, $lastPara.Range.End)\n        $delRange.Delete()\n    }\n}\n"}
